$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (row 1) — new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, border, centered) from AC1 onto the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows 2-44 — team record values
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 85   # AD
    $ws.Cells.Item($r, 31).Value = 77   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
